# Actualización automática del index.html y archivo Excel
#
# The INCO claims sheet dropped two entries that are no longer part of the
# published feed:
#   - Row 15: Caso -159, "WARNES ,AV. /ALT/ 1605" (Chacarita)
#   - Row 36: Caso -252, "LIBERTI TOMAS /ALT/ 1110" (La Boca)
#
# Removing these two rows shifts every following row up by one position
# (twice), which is exactly what the diff shows (e.g. old row 16 data now
# lives in row 15, old row 37 data now lives in row 35, etc.), and shrinks
# the used range from A1:N42 down to A1:N40.
#
# Delete the higher-numbered row first so the row index of the still-to-be
# -deleted row (15) doesn't shift before we get to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(36).EntireRow.Delete()
$ws.Rows(15).EntireRow.Delete()
